# Add 'avenant revision loyer EA' functionality
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: renew contract AV4 -> AV6, with updated avenant amounts ---
$ws.Range("G2").Value = "044/FES VILLE /AV6"
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 9400

# --- Row 3: fill in the previously blank row with a second avenant line ---
$ws.Range("A3").Value = "KHADIJA LALA"
$ws.Range("B3").Value = "K5443645"
$ws.Range("C3").Value = "'354564564324158786713544"
$ws.Range("D3").Value = "AG 100"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "044/FES VILLE /AV6"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 30000
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 25500

# --- Row 4: new third avenant line ---
$ws.Range("A4").Value = "KHADIJA LALA"
$ws.Range("B4").Value = "K5443645"
$ws.Range("C4").Value = "'354564564324158786713544"
$ws.Range("D4").Value = "AG 100"
$ws.Range("E4").Value = "BP"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "044/FES VILLE /AV6"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 4300

# --- Row 5: totals row (blank text columns, summed amount columns) ---
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("I5").Value = 45000
$ws.Range("J5").Value = 5800
$ws.Range("K5").Value = 39200
